# Populate the employee test row and restore the sheet to an
# "unprotected, default view" state, matching the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell data (row 2): A-C "test", D blank, E "employee_5", F a number ---
$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = "test"
$ws.Range("C2").Value = "test"
$ws.Range("E2").Value = "employee_5"
$ws.Range("F2").Value = 4345528

# --- Sheet-level settings / "additional settings" ---

# Remove sheet protection entirely (mark system / unrestricted editing).
$ws.Unprotect()

# Outline defaults: summary rows below, summary columns to the right.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Reset selection back to A1 (top-left), like a freshly opened sheet.
$ws.Range("A1").Select() | Out-Null

# Page setup basics.
$ws.PageSetup.PaperSize = 1
$ws.PageSetup.Zoom = 100
